# results_dota.xlsx update
#   "Added more results, added class sampling to chipper"
#
# Concrete, scriptable parts of the change (window-geometry / revision
# GUIDs in xl/workbook.xml are Excel session bookkeeping that isn't
# exposed for writing through this object model, so they're left alone):
#
#   1. Sheet1: columns B:I are resized from their old "AutoFit"/bestFit
#      widths down to a narrower, explicit width (bestFit is cleared as
#      a side effect of assigning ColumnWidth directly).
#   2. Chart "Train/Test Accuracy (DOTA)" (the first chart object on
#      Sheet1): the value (Y) axis minimum scale moves from 40 to 50.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Narrow columns B:I ------------------------------------------------
# Target stored widths (characters) from the authored workbook:
#   B=12.28515625 C=13.140625 D=12.7109375 E=12.85546875
#   F=13.7109375  G=11.28515625 H=13.42578125 I=13.140625
$ws.Columns.Item(2).ColumnWidth = 11.5
$ws.Columns.Item(3).ColumnWidth = 12.33333333333333
$ws.Columns.Item(4).ColumnWidth = 11.83333333333333
$ws.Columns.Item(5).ColumnWidth = 12
$ws.Columns.Item(6).ColumnWidth = 12.83333333333333
$ws.Columns.Item(7).ColumnWidth = 10.5
$ws.Columns.Item(8).ColumnWidth = 12.66666666666667
$ws.Columns.Item(9).ColumnWidth = 12.33333333333333

# --- 2. Rescale the accuracy chart's value axis --------------------------
$chartObj = $ws.ChartObjects(1)
$chart = $chartObj.Chart
$valueAxis = $chart.Axes(2)          # xlValue axis (the category axis is Axes(1))
$valueAxis.MinimumScale = 50
